# Auto-generated edit script applying numeric updates to Halicarnassus_Profits sheets
# (ALC, ARM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2651.2188
$ws.Range("I17").Value = 1500
$ws.Range("J17").Value = 2770.3103
$ws.Range("K17").Value = 4500
$ws.Range("L17").Value = 8310.930899999999
$ws.Range("M17").Value = -4332
$ws.Range("N17").Value = -8646.930899999999

# Row 70
$ws.Range("H70").Value = 7660.222
$ws.Range("I70").Value = 1999.75
$ws.Range("K70").Value = 5999.25
$ws.Range("M70").Value = -5729.25

# Row 73
$ws.Range("H73").Value = 7660.222
$ws.Range("I73").Value = 1999.75
$ws.Range("K73").Value = 5999.25
$ws.Range("M73").Value = -5063.25

# Row 97
$ws.Range("H97").Value = 790.3333
$ws.Range("J97").Value = 790.3333
$ws.Range("L97").Value = 2370.9999
$ws.Range("N97").Value = -3362.9999

# Row 100
$ws.Range("H100").Value = 2020.375
$ws.Range("I100").Value = 803.1
$ws.Range("J100").Value = 4049.1667
$ws.Range("K100").Value = 803.1
$ws.Range("L100").Value = 4049.1667
$ws.Range("M100").Value = -262.1
$ws.Range("N100").Value = -5131.1667

# Row 104
$ws.Range("H104").Value = 1013.4286
$ws.Range("I104").Value = 1015.6667
$ws.Range("K104").Value = 3047.0001
$ws.Range("M104").Value = -1300.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 122
$ws.Range("H122").Value = 2771.5557
$ws.Range("I122").Value = 2680.5
$ws.Range("K122").Value = 8041.5
$ws.Range("M122").Value = -5591.5

$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 1174.5555
$ws.Range("I5").Value = 262.16666
$ws.Range("K5").Value = 262.16666
$ws.Range("M5").Value = -150.16666

# Row 7
$ws.Range("H7").Value = 59
$ws.Range("I7").Value = 58.77778
$ws.Range("K7").Value = 58.77778
$ws.Range("M7").Value = 54.22222

# Row 68
$ws.Range("H68").Value = 79953
$ws.Range("J68").Value = 84941.25
$ws.Range("L68").Value = 84941.25
$ws.Range("N68").Value = -86439.25

# Row 71
$ws.Range("H71").Value = 79953
$ws.Range("J71").Value = 84941.25
$ws.Range("L71").Value = 254823.75
$ws.Range("N71").Value = -262311.75

# Row 86
$ws.Range("H86").Value = 3584
$ws.Range("I86").Value = 3440.8
$ws.Range("K86").Value = 3440.8
$ws.Range("M86").Value = -2317.8

# Row 89
$ws.Range("H89").Value = 3584
$ws.Range("I89").Value = 3440.8
$ws.Range("K89").Value = 17204
$ws.Range("M89").Value = -11588

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 141
$ws.Range("I14").Value = 141
$ws.Range("K14").Value = 423
$ws.Range("M14").Value = -250

# Row 70
$ws.Range("H70").Value = 93
$ws.Range("I70").Value = 93
$ws.Range("K70").Value = 279
$ws.Range("M70").Value = 36

# Row 73
$ws.Range("H73").Value = 93
$ws.Range("I73").Value = 93
$ws.Range("K73").Value = 279
$ws.Range("M73").Value = 813

# Row 80
$ws.Range("H80").Value = 4781.4
$ws.Range("I80").Value = 4784.95
$ws.Range("J80").Value = 4767.2
$ws.Range("K80").Value = 14354.85
$ws.Range("L80").Value = 14301.6
$ws.Range("M80").Value = -13418.85
$ws.Range("N80").Value = -16173.6

# Row 83
$ws.Range("H83").Value = 4781.4
$ws.Range("I83").Value = 4784.95
$ws.Range("J83").Value = 4767.2
$ws.Range("K83").Value = 43064.55
$ws.Range("L83").Value = 42904.8
$ws.Range("M83").Value = -38384.55
$ws.Range("N83").Value = -52264.8

# Row 107
$ws.Range("H107").Value = 754.4286
$ws.Range("I107").Value = 491
$ws.Range("J107").Value = 859.8
$ws.Range("K107").Value = 1473
$ws.Range("L107").Value = 2579.4
$ws.Range("M107").Value = 447
$ws.Range("N107").Value = -6419.4

# Row 132
$ws.Range("H132").Value = 1309.6666
$ws.Range("I132").Value = 972.25
$ws.Range("K132").Value = 8750.25
$ws.Range("M132").Value = -6220.25

# Row 139
$ws.Range("H139").Value = 5969.6
$ws.Range("I139").Value = 5969.6
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 17908.8
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -12768.8
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 12917117
$ws.Range("I11").Value = 10625450
$ws.Range("J11").Value = 17500450
$ws.Range("K11").Value = 10625450
$ws.Range("L11").Value = 17500450
$ws.Range("M11").Value = -10625311
$ws.Range("N11").Value = -17500728

# Row 18
$ws.Range("H18").Value = 100000
$ws.Range("J18").Value = 100000
$ws.Range("L18").Value = 100000
$ws.Range("N18").Value = -100586

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 607.5714
$ws.Range("I16").Value = 375.5
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 375.5
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -205.5
$ws.Range("N16").Value = -2340

# Row 20
$ws.Range("H20").Value = 999999
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 999999
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 999999
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -1000451

# Row 22
$ws.Range("H22").Value = 812.61536
$ws.Range("J22").Value = 760.75
$ws.Range("L22").Value = 760.75
$ws.Range("N22").Value = -1350.75

# Row 27
$ws.Range("H27").Value = 812.61536
$ws.Range("J27").Value = 760.75
$ws.Range("L27").Value = 760.75
$ws.Range("N27").Value = -974.75

# Row 40
$ws.Range("H40").Value = 7853.6665
$ws.Range("I40").Value = 7853.6665
$ws.Range("K40").Value = 7853.6665
$ws.Range("M40").Value = -7717.6665

$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 17250.75
$ws.Range("I3").Value = 4999.5
$ws.Range("J3").Value = 29502
$ws.Range("K3").Value = 4999.5
$ws.Range("L3").Value = 29502
$ws.Range("M3").Value = -4885.5
$ws.Range("N3").Value = -29730

# Row 10
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

# Row 20
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# Row 122
$ws.Range("H122").Value = 3281.8
$ws.Range("I122").Value = 1634.6666
$ws.Range("K122").Value = 4903.9998
$ws.Range("M122").Value = -2453.9998

# Row 126
$ws.Range("H126").Value = 3568.261
$ws.Range("I126").Value = 1497.9286
$ws.Range("K126").Value = 4493.7858
$ws.Range("M126").Value = -2023.7858
